$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2333.2222
$ws.Range("I43").Value = 1999
$ws.Range("J43").Value = 2375
$ws.Range("K43").Value = 1999
$ws.Range("L43").Value = 2375
$ws.Range("M43").Value = -1930
$ws.Range("N43").Value = -2513
$ws.Range("H92").Value = 17127.25
$ws.Range("I92").Value = 24250
$ws.Range("J92").Value = 10004.5
$ws.Range("K92").Value = 24250
$ws.Range("L92").Value = 10004.5
$ws.Range("M92").Value = -23002
$ws.Range("N92").Value = -12500.5
$ws.Range("H112").Value = 7648.4473
$ws.Range("J112").Value = 10505.423
$ws.Range("L112").Value = 31516.269
$ws.Range("N112").Value = -33732.269
$ws.Range("H116").Value = 333267.53
$ws.Range("I116").Value = 9000
$ws.Range("J116").Value = 454867.88
$ws.Range("K116").Value = 9000
$ws.Range("L116").Value = 454867.88
$ws.Range("M116").Value = -5558
$ws.Range("N116").Value = -461751.88
$ws.Range("H132").Value = 53448.914
$ws.Range("I132").Value = 66944
$ws.Range("K132").Value = 200832
$ws.Range("M132").Value = -198302
$ws.Range("H137").Value = 1958.6364
$ws.Range("I137").Value = 1596.1786
$ws.Range("K137").Value = 4788.5358
$ws.Range("M137").Value = -2238.5358
$ws.Range("H138").Value = 3067.1904
$ws.Range("I138").Value = 3696.25
$ws.Range("J138").Value = 2919.1765
$ws.Range("K138").Value = 11088.75
$ws.Range("L138").Value = 8757.529500000001
$ws.Range("M138").Value = -5948.75
$ws.Range("N138").Value = -19037.5295

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9618770
$ws.Range("I32").Value = 10871875
$ws.Range("J32").Value = 11630
$ws.Range("K32").Value = 10871875
$ws.Range("L32").Value = 11630
$ws.Range("M32").Value = -10871588
$ws.Range("N32").Value = -12204
$ws.Range("H62").Value = 34547.5
$ws.Range("J62").Value = 34547.5
$ws.Range("L62").Value = 34547.5
$ws.Range("N62").Value = -35795.5
$ws.Range("H65").Value = 34547.5
$ws.Range("J65").Value = 34547.5
$ws.Range("L65").Value = 103642.5
$ws.Range("N65").Value = -109882.5
$ws.Range("H74").Value = 3380547.2
$ws.Range("I74").Value = 4311782.5
$ws.Range("J74").Value = 4819.5
$ws.Range("K74").Value = 4311782.5
$ws.Range("L74").Value = 4819.5
$ws.Range("M74").Value = -4310908.5
$ws.Range("N74").Value = -6567.5
$ws.Range("H77").Value = 3380547.2
$ws.Range("I77").Value = 4311782.5
$ws.Range("J77").Value = 4819.5
$ws.Range("K77").Value = 21558912.5
$ws.Range("L77").Value = 24097.5
$ws.Range("M77").Value = -21554544.5
$ws.Range("N77").Value = -32833.5
$ws.Range("H122").Value = 2021.0851
$ws.Range("I122").Value = 1922.6136
$ws.Range("K122").Value = 5767.8408
$ws.Range("M122").Value = -3317.8408
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 41999.75
$ws.Range("J51").Value = 41999.75
$ws.Range("L51").Value = 41999.75
$ws.Range("N51").Value = -42981.75
$ws.Range("H55").Value = 77999.5
$ws.Range("J55").Value = 77999.5
$ws.Range("L55").Value = 77999.5
$ws.Range("N55").Value = -78545.5
$ws.Range("H94").Value = 1142.7727
$ws.Range("I94").Value = 952.5
$ws.Range("J94").Value = 1999
$ws.Range("K94").Value = 952.5
$ws.Range("L94").Value = 1999
$ws.Range("M94").Value = -501.5
$ws.Range("N94").Value = -2901
$ws.Range("H105").Value = 1285.2778
$ws.Range("I105").Value = 1411.1428
$ws.Range("K105").Value = 1411.1428
$ws.Range("M105").Value = 335.8571999999999
$ws.Range("H134").Value = 1328384.4
$ws.Range("I134").Value = 2980316
$ws.Range("K134").Value = 8940948
$ws.Range("M134").Value = -8938413

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 67217.46000000001
$ws.Range("I31").Value = 105823.39
$ws.Range("J31").Value = 19346.12
$ws.Range("K31").Value = 105823.39
$ws.Range("L31").Value = 19346.12
$ws.Range("M31").Value = -105528.39
$ws.Range("N31").Value = -19936.12
$ws.Range("H34").Value = 67217.46000000001
$ws.Range("I34").Value = 105823.39
$ws.Range("J34").Value = 19346.12
$ws.Range("K34").Value = 105823.39
$ws.Range("L34").Value = 19346.12
$ws.Range("M34").Value = -105621.39
$ws.Range("N34").Value = -19750.12
$ws.Range("H99").Value = 6392.3335
$ws.Range("I99").Value = 6613.5
$ws.Range("J99").Value = 5950
$ws.Range("K99").Value = 6613.5
$ws.Range("L99").Value = 5950
$ws.Range("M99").Value = -5115.5
$ws.Range("N99").Value = -8946
$ws.Range("H126").Value = 6392.3335
$ws.Range("I126").Value = 6613.5
$ws.Range("J126").Value = 5950
$ws.Range("K126").Value = 19840.5
$ws.Range("L126").Value = 17850
$ws.Range("M126").Value = -17370.5
$ws.Range("N126").Value = -22790
$ws.Range("H129").Value = 34868.75
$ws.Range("I129").Value = 37245
$ws.Range("J129").Value = 32492.5
$ws.Range("K129").Value = 37245
$ws.Range("L129").Value = 32492.5
$ws.Range("M129").Value = -32245
$ws.Range("N129").Value = -42492.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 37.411766
$ws.Range("I2").Value = 30.214285
$ws.Range("K2").Value = 181.28571
$ws.Range("M2").Value = -68.28570999999999
$ws.Range("H75").Value = 6040.9165
$ws.Range("I75").Value = 550
$ws.Range("K75").Value = 1650
$ws.Range("M75").Value = -652
$ws.Range("H78").Value = 6040.9165
$ws.Range("I78").Value = 550
$ws.Range("K78").Value = 4950
$ws.Range("M78").Value = 42
$ws.Range("H99").Value = 5561.5835
$ws.Range("I99").Value = 3328.375
$ws.Range("J99").Value = 10028
$ws.Range("K99").Value = 9985.125
$ws.Range("L99").Value = 30084
$ws.Range("M99").Value = -7739.125
$ws.Range("N99").Value = -34576
$ws.Range("H113").Value = 672.4
$ws.Range("I113").Value = 385.8
$ws.Range("J113").Value = 959
$ws.Range("K113").Value = 1157.4
$ws.Range("L113").Value = 2877
$ws.Range("M113").Value = 1012.6
$ws.Range("N113").Value = -7217
$ws.Range("H118").Value = 3048.5
$ws.Range("I118").Value = 1733
$ws.Range("J118").Value = 6995
$ws.Range("K118").Value = 5199
$ws.Range("L118").Value = 20985
$ws.Range("M118").Value = -3956
$ws.Range("N118").Value = -23471

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 39999.5
$ws.Range("J26").Value = 39999.5
$ws.Range("L26").Value = 39999.5
$ws.Range("N26").Value = -40559.5
$ws.Range("H50").Value = 39999.5
$ws.Range("J50").Value = 39999.5
$ws.Range("L50").Value = 39999.5
$ws.Range("N50").Value = -40995.5
$ws.Range("H52").Value = 50448.875
$ws.Range("J52").Value = 50448.875
$ws.Range("L52").Value = 50448.875
$ws.Range("N52").Value = -50966.875
$ws.Range("H97").Value = 3474.5
$ws.Range("I97").Value = 3179.3333
$ws.Range("J97").Value = 4360
$ws.Range("K97").Value = 3179.3333
$ws.Range("L97").Value = 4360
$ws.Range("M97").Value = -2683.3333
$ws.Range("N97").Value = -5352
$ws.Range("H102").Value = 1896.3914
$ws.Range("I102").Value = 1442.8823
$ws.Range("K102").Value = 1442.8823
$ws.Range("M102").Value = 179.1177
$ws.Range("H123").Value = 33999.5
$ws.Range("J123").Value = 33999.5
$ws.Range("L123").Value = 33999.5
$ws.Range("N123").Value = -38899.5
$ws.Range("H126").Value = 620793.75
$ws.Range("I126").Value = 728145.8
$ws.Range("K126").Value = 2184437.4
$ws.Range("M126").Value = -2181967.4
$ws.Range("H131").Value = 52500
$ws.Range("J131").Value = 52500
$ws.Range("L131").Value = 52500
$ws.Range("N131").Value = -62580

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2229.6
$ws.Range("I22").Value = 1999.5
$ws.Range("J22").Value = 3150
$ws.Range("K22").Value = 1999.5
$ws.Range("L22").Value = 3150
$ws.Range("M22").Value = -1704.5
$ws.Range("N22").Value = -3740
$ws.Range("H27").Value = 2229.6
$ws.Range("I27").Value = 1999.5
$ws.Range("J27").Value = 3150
$ws.Range("K27").Value = 1999.5
$ws.Range("L27").Value = 3150
$ws.Range("M27").Value = -1892.5
$ws.Range("N27").Value = -3364
$ws.Range("H35").Value = 1906.5264
$ws.Range("J35").Value = 1199
$ws.Range("L35").Value = 1199
$ws.Range("N35").Value = -1871
$ws.Range("H40").Value = 3109.2727
$ws.Range("I40").Value = 2880.8462
$ws.Range("K40").Value = 2880.8462
$ws.Range("M40").Value = -2744.8462
$ws.Range("H108").Value = 72266.125
$ws.Range("J108").Value = 72266.125
$ws.Range("L108").Value = 72266.125
$ws.Range("N108").Value = -79946.125
$ws.Range("H122").Value = 3352.05
$ws.Range("J122").Value = 3878.4285
$ws.Range("L122").Value = 11635.2855
$ws.Range("N122").Value = -16535.2855
$ws.Range("H132").Value = 845654.3
$ws.Range("I132").Value = 1194525.5
$ws.Range("K132").Value = 3583576.5
$ws.Range("M132").Value = -3581046.5
$ws.Range("H136").Value = 51883.96
$ws.Range("I136").Value = 2744.5
$ws.Range("K136").Value = 8233.5
$ws.Range("M136").Value = -5683.5

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 73449.75
$ws.Range("J16").Value = 73449.75
$ws.Range("L16").Value = 73449.75
$ws.Range("N16").Value = -74033.75
$ws.Range("H126").Value = 2994
$ws.Range("I126").Value = 3368
$ws.Range("J126").Value = 1622.6666
$ws.Range("K126").Value = 10104
$ws.Range("L126").Value = 4867.9998
$ws.Range("M126").Value = -7634
$ws.Range("N126").Value = -9807.9998
$ws.Range("H132").Value = 3050796.2
$ws.Range("I132").Value = 3146090.5
$ws.Range("K132").Value = 9438271.5
$ws.Range("M132").Value = -9435741.5
$ws.Range("H136").Value = 5853179.5
$ws.Range("I136").Value = 6910397
$ws.Range("J136").Value = 38484
$ws.Range("K136").Value = 20731191
$ws.Range("L136").Value = 115452
$ws.Range("M136").Value = -20728641
